$wb = $excel.ActiveWorkbook

# Sheet "2025" - row 2 updates
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.003725587022858967
$ws.Range("E2").Value = 0.3717629768396423
$ws.Range("G2").Value = 0.2494892361375005
$ws.Range("I2").Value = 0.3678009020113525
$ws.Range("L2").Value = 0.5957133999999999
$ws.Range("M2").Value = 0.08226208333333333
$ws.Range("N2").Value = 12.82873417801238
$ws.Range("O2").Value = 3.475425528524659

# Sheet "2030" - row 2 updates
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.04671422655046342
$ws.Range("E2").Value = 0.227037787022859
$ws.Range("I2").Value = 0.5575521882485598
$ws.Range("L2").Value = 0.1504622697400875
$ws.Range("M2").Value = 0.04790916666666666
$ws.Range("N2").Value = 5.725818980985812
$ws.Range("O2").Value = 2.958482383999864

# Sheet "2035" - row 2 updates
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.09246557211163992
$ws.Range("B2").Value = 0.03201577550183085
$ws.Range("E2").Value = 0.1659256519093786
$ws.Range("I2").Value = 0.2978005395653942
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04688944789627539
$ws.Range("N2").Value = 7.70534674480178
$ws.Range("O2").Value = 4.371478642922033

$wb.Save()
